$d = $word.ActiveDocument

# NOTE: all paragraph indices below refer to the ORIGINAL (pre-edit) document.
# Operations are applied from the bottom of the document upward so that the
# indices used for earlier (lower-numbered) paragraphs are never invalidated
# by edits made further down.

# ---------------------------------------------------------------------------
# Occurrences section (tail of the list)
# ---------------------------------------------------------------------------

# Final paragraph: repurpose the "(ITransformResource, ITransformKind, ...)" shell
$d.Paragraphs(142).Range.Text = "Contexts Resources, Kinds, Occurrences: Statement (relation data), Mapping (schema), Transform (behavior) given composite SK(PK, OK), PK(SK, OK), OK(PK, SK)."

# Remove IStatementOccurrence/IMappingOccurrence/ITransformOccurrence quadruple (5 paragraphs)
$rng = $d.Range($d.Paragraphs(137).Range.Start, $d.Paragraphs(141).Range.End)
$rng.Delete()

$d.Paragraphs(135).Range.Text = "IObjectOccurrence : IObject"
$d.Paragraphs(133).Range.Text = "IPredicateOccurrence : IPredicate"
$d.Paragraphs(132).Range.Text = "(ISubjectResource, ISubjectKind, IPredicateOccurrence, IObjectOccurrence) templates hierarchy in / out statements"
$d.Paragraphs(131).Range.Text = "ISubjectOccurrence : ISubject, Context<Subject> / Subject<Context>"

$d.Paragraphs(130).Range.Text = "Occurrences (Contexts):"
$d.Paragraphs(130).Range.InsertParagraphAfter()
$d.Paragraphs(131).Range.Text = "IContextOccurrence : IContext"
$d.Paragraphs(131).Range.InsertParagraphAfter()
$d.Paragraphs(132).Range.Text = "(IContextResource, IContextKind, IPredicateKind, IObjectKind)"

# ---------------------------------------------------------------------------
# Kinds section
# ---------------------------------------------------------------------------

# Remove IStatementKind/IMappingKind/ITransformKind quadruple + trailing blank (7 paragraphs)
$rng = $d.Range($d.Paragraphs(123).Range.Start, $d.Paragraphs(129).Range.End)
$rng.Delete()

$d.Paragraphs(118).Range.Text = "(IContextOccurrence, SK of PK/OK: Relation, PK, OK)"
$d.Paragraphs(118).Range.InsertParagraphAfter()
$d.Paragraphs(119).Range.Text = "ISubjectKind : Kind<Subject> / Subject<Kind>"
$d.Paragraphs(119).Range.InsertParagraphAfter()
$d.Paragraphs(120).Range.Text = "(ISubjectOccurrence, ISubjectKind, IPredicateResource, IObjectResource) templates hierarchy in / out statements"

$d.Paragraphs(117).Range.Text = "IContextKind"

# ---------------------------------------------------------------------------
# Resources section
# ---------------------------------------------------------------------------

# Remove IStatementResource/IMappingResource/ITransformResource quadruple (6 paragraphs)
$rng = $d.Range($d.Paragraphs(109).Range.Start, $d.Paragraphs(114).Range.End)
$rng.Delete()

$d.Paragraphs(107).Range.Text = "IObjectResource : IObject"
$d.Paragraphs(105).Range.Text = "IPredicateResource : IPredicate"

$d.Paragraphs(104).Range.Text = "(IContextKind, IContextOccurrence, IPredicateKind, IObjectKind)"
$d.Paragraphs(104).Range.InsertParagraphAfter()
$d.Paragraphs(105).Range.Text = "ISubjectResource : ISubject, Subject<Resource> / Resource<Subject>"
$d.Paragraphs(105).Range.InsertParagraphAfter()
$d.Paragraphs(106).Range.Text = "(ISubjectKind, ISubjectOccurrence, IPredicateKind, IObjectKind) templates hierarchy in / out statements"

$d.Paragraphs(103).Range.Text = "IContextResource : IContext"

# ---------------------------------------------------------------------------
# Interfaces (Sets) section
# ---------------------------------------------------------------------------

$d.Paragraphs(94).Range.Text = "IContext : ISubject, IPredicate, IObject, ISubjectKind, IObjectKind, IPredicateKind"
